# Updates the "cryptos" price/volume table (columns D and E) with refreshed
# market data, and swaps the ARBITRUM / MXToken rows (48-49).
#
# Note: several Price values (column D) are plain decimal-looking strings
# (e.g. "227.13") that must stay as literal text, matching the source
# workbook where these are text cells, not numbers. Prefixing the value
# with a leading apostrophe (doubled to '' inside a single-quoted
# PowerShell string) forces Excel to store it as text instead of
# auto-converting it to a floating point number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.711.57'
$ws.Range('E2').Value = '  -1.46%  '

$ws.Range('D3').Value = '2.025.45'
$ws.Range('E3').Value = '  -1.90%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').Value = '''227.13'
$ws.Range('E5').Value = '  -1.51%  '

$ws.Range('D6').Value = '''0.605'
$ws.Range('E6').Value = '  -1.89%  '

$ws.Range('D7').Value = '''59.97'
$ws.Range('E7').Value = '  -1.78%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('E9').Value = '  -3.68%  '

$ws.Range('D10').Value = '''0.0820'
$ws.Range('E10').Value = '  +0.87%  '

$ws.Range('D11').Value = '''0.103'
$ws.Range('E11').Value = '  -0.81%  '

$ws.Range('D12').Value = '2.325.33'
$ws.Range('E12').Value = '  -1.75%  '

$ws.Range('D13').Value = '''14.37'
$ws.Range('E13').Value = '  -3.28%  '

$ws.Range('D14').Value = '''20.99'
$ws.Range('E14').Value = '  -2.13%  '

$ws.Range('D15').Value = '''0.757'
$ws.Range('E15').Value = '  -0.99%  '

$ws.Range('D16').Value = '''5.16'
$ws.Range('E16').Value = '  -3.33%  '

$ws.Range('D17').Value = '2.022.26'
$ws.Range('E17').Value = '  -2.06%  '

$ws.Range('D18').Value = '37.627.73'
$ws.Range('E18').Value = '  -1.43%  '

$ws.Range('D19').Value = '''69.40'
$ws.Range('E19').Value = '  -1.12%  '

$ws.Range('D20').Value = '''5.90'
$ws.Range('E20').Value = '  -6.68%  '

$ws.Range('D21').Value = '0.0₃0822'
$ws.Range('E21').Value = '  -1.83%  '

$ws.Range('D22').Value = '''223.36'
$ws.Range('E22').Value = '  -1.19%  '

$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('D24').Value = '''2.38'
$ws.Range('E24').Value = '  -1.94%  '

$ws.Range('D25').Value = '''2.25'
$ws.Range('E25').Value = '  -0.01%  '

$ws.Range('D26').Value = '''167.36'
$ws.Range('E26').Value = '  +0.70%  '

$ws.Range('D27').Value = '''9.25'
$ws.Range('E27').Value = '  -0.95%  '

$ws.Range('E28').Value = '  -3.55%  '

$ws.Range('D29').Value = '''18.77'
$ws.Range('E29').Value = '  -1.99%  '

$ws.Range('E30').Value = '  -5.68%  '

$ws.Range('E31').Value = '  -0.28%  '

$ws.Range('D32').Value = '''2.20'
$ws.Range('E32').Value = '  +7.50%  '

$ws.Range('D33').Value = '''4.38'
$ws.Range('E33').Value = '  -4.23%  '

$ws.Range('D34').Value = '''0.0602'
$ws.Range('E34').Value = '  -1.00%  '

$ws.Range('D35').Value = '''4.46'
$ws.Range('E35').Value = '  -3.93%  '

$ws.Range('D36').Value = '''6.36'
$ws.Range('E36').Value = '  +1.97%  '

$ws.Range('D37').Value = '''2.30'
$ws.Range('E37').Value = '  -1.16%  '

$ws.Range('D38').Value = '''3.36'
$ws.Range('E38').Value = '  +0.79%  '

$ws.Range('E39').Value = '  +0.04%  '

$ws.Range('D40').Value = '''17.92'
$ws.Range('E40').Value = '  +4.27%  '

$ws.Range('D41').Value = '1.534.19'
$ws.Range('E41').Value = '  +0.29%  '

$ws.Range('E42').Value = '  -1.31%  '

$ws.Range('D43').Value = '''95.42'
$ws.Range('E43').Value = '  -2.97%  '

$ws.Range('D44').Value = '''2.81'
$ws.Range('E44').Value = '  -2.31%  '

$ws.Range('E45').Value = '  -2.27%  '

$ws.Range('D46').Value = '''4.08'
$ws.Range('E46').Value = '  +0.22%  '

$ws.Range('E47').Value = '  -3.07%  '

$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = '''1.00'
$ws.Range('E48').Value = '  -2.39%  '

$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').Value = '''2.96'
$ws.Range('E49').Value = '  -0.01%  '

$ws.Range('D50').Value = '''7.07'
$ws.Range('E50').Value = '  -0.75%  '

$ws.Range('D51').Value = '2.215.76'
$ws.Range('E51').Value = '  -1.72%  '
